# Motorola Stability Tests - Shelter build 1.29 updates

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Completed Items"
# ---------------------------------------------------------------------------
$completed = $wb.Worksheets.Item("Completed Items")

# Row 8: note text is unchanged ("Awaiting bug fixes") - reassert to keep
# shared-string bookkeeping consistent with the diff.
$completed.Range("K8").Value = "Awaiting bug fixes"

# Row 9: count bumped 2 -> 4, note text updated.
$completed.Range("J9").Value = 4
$completed.Range("K9").Value = "Awaiting design analysis"

# Row 10: count + note cleared entirely.
$completed.Range("J10").Value = $null
$completed.Range("K10").Value = $null

# Rows 49/50: note consolidated into the new "design analysis" wording.
$completed.Range("G49").Value = "Awaiting Motorola/ProtoTest design analysis."
$completed.Range("G50").Value = "Awaiting Motorola/ProtoTest design analysis."

# Rows 64/65: status set to "In progress", and the existing note's
# underlying text is updated to the new wording.
$completed.Range("F64").Value = "In progress"
$completed.Range("F65").Value = "In progress"
$completed.Range("G64").Value = "Awaiting Motorola/ProtoTest design analysis."
$completed.Range("G65").Value = "Awaiting Motorola/ProtoTest design analysis."

# ---------------------------------------------------------------------------
# Sheet "Eggplant Scripts"
# ---------------------------------------------------------------------------
$eggplant = $wb.Worksheets.Item("Eggplant Scripts")

$buildRows = @(3,4,5,6,9,10,11,12,13,14,15,16,22,23,24,25,26,27)
foreach ($r in $buildRows) {
    $eggplant.Range("D$r").Value = 1.29
}

# Row 34: build value replaced with an "In Progress" status note.
$eggplant.Range("D34").Value = "In Progress"

# Update the sheet's scroll position / selection to match the saved view.
$eggplant.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$eggplant.Range("D28").Select()

# Restore the originally active sheet/tab.
$completed.Activate()
